# Apply the "Project link" slide update:
#  - Reposition + restyle the title ("Project link" -> "Project link :",
#    single underline, Sitka Display font)
#  - Move the body placeholder
#  - Insert the hyperlinked GitHub project URL into the body placeholder
#
# NOTE on the literal position/size numbers below: PowerPoint's COM object
# model stores Shape.Left/Top/Width/Height as single-precision (float32)
# points, so the EMU value that ultimately lands in the XML is
# int(float64(float32(points)) * 12700). The decimal literals used here were
# chosen so that round-trip lands exactly on the target EMU values from the
# diff (520979/1124744/10681335/758190 for the title, 520979/2708920 for the
# body offset).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

$title = $s.Shapes.Item(1)
$body  = $s.Shapes.Item(2)

# --- Title shape -----------------------------------------------------
$title.Left   = 41.021968841552734
$title.Top    = 88.5625228881836
$title.Width  = 841.050048828125
$title.Height = 59.70000076293945

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Project link :"
$titleRange.Font.Underline = $true
$titleRange.Font.Name = "Sitka Display"

# --- Body placeholder shape -------------------------------------------
$body.Left = 41.021968841552734
$body.Top  = 213.30079650878906

$bodyRange = $body.TextFrame.TextRange
$linkText = "https://github.com/PrabhakarPULIGADDA/prabhakar-Project.git"
$newRun = $bodyRange.InsertBefore($linkText)
$newRun.LanguageID = "en-US"
$newRun.Font.Size = 28
$newRun.ActionSettings(1).Hyperlink.Address = $linkText
